$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new, otherwise unused, sheet "Feuille inutile" between
#    "2005" and "2006".
#    NOTE: worksheet references are index-bound in this runtime, so any
#    sheet handle captured before an Add/rename that shifts tab order must
#    be re-fetched (by name) afterwards instead of reused.
# ---------------------------------------------------------------------------
$sheet2005 = $wb.Worksheets.Item("2005")

$feuilleInutile = $wb.Worksheets.Add($null, $sheet2005)
$feuilleInutile.Name = "Feuille inutile"

# Re-fetch "2006" now that the tab order has shifted.
$sheet2006 = $wb.Worksheets.Item("2006")

$leftStyleSrc  = $sheet2006.Range("A1")   # style index 1 in the original file (left aligned)
$rightStyleSrc = $sheet2006.Range("B2")   # style index 2 in the original file (right aligned)

# ---------------------------------------------------------------------------
# 2. Rebuild sheet "2005" with the new header/title rows followed by the
#    pre-existing data rows (now shifted down to rows 5-8).
# ---------------------------------------------------------------------------
$sheet2005.Cells.Clear()

# Row 1-2 : plain title rows ("Automatic" text colour, Arial, bottom aligned)
$titleRange = $sheet2005.Range("A1:B2")
$titleRange.Font.Name = "Arial"
$titleRange.Font.ThemeColor = 1

$sheet2005.Range("A1").Value = "Parcours thématique_RGPD & Données personnelles"
$sheet2005.Range("A2").Value = "Création : Avril 2022"

# Row 3 : dark-grey filled banner row, white Arial text
$bannerRange = $sheet2005.Range("A3:B3")
$bannerRange.Interior.Pattern = 1
$bannerRange.Interior.PatternColor = 6710886
$bannerRange.Interior.Color = 6710886
$bannerRange.Font.Color = 16777215
$bannerRange.Font.Name = "Arial"

$sheet2005.Range("A3").Value = "1) capage"

# Row 4 : bold column headers
$headerRange = $sheet2005.Range("A4:B4")
$headerRange.Font.Bold = $true
$headerRange.Font.ThemeColor = 1

$sheet2005.Range("A4").Value = "Sujet (code ou autre)"
$sheet2005.Range("B4").Value = "Niveau max"

# Rows 5-8 : original data rows, re-using the two pre-existing cell styles.
$sheet2005.Range("A5").Value = "@codageEmblématique"
$sheet2005.Range("A6").Value = "@terminal"
$sheet2005.Range("A7").Value = "@editerDocEnLigne"
$sheet2005.Range("A8").Value = "@partageDroits"

$sheet2005.Range("B5").Value = 1
$sheet2005.Range("B6").Value = 2
$sheet2005.Range("B7").Value = 3
$sheet2005.Range("B8").Value = 1

$leftStyleSrc.Copy()
$sheet2005.Range("A5:A8").PasteSpecial(-4122)

$rightStyleSrc.Copy()
$sheet2005.Range("B5:B8").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Populate the new "Feuille inutile" sheet.
# ---------------------------------------------------------------------------
$funRange = $feuilleInutile.Range("A1:B1")
$funRange.Font.Name = "Arial"
$funRange.Font.ThemeColor = 1

$feuilleInutile.Range("A1").Value = "Coucou"
$feuilleInutile.Range("B1").Value = "les"

$funRange.Copy()
$feuilleInutile.Range("A2").PasteSpecial(-4122)
$feuilleInutile.Range("B2").Clear()
$feuilleInutile.Range("A2").Value = "zamiiiiiiiiiiiiis"

$excel.CutCopyMode = $false
